$wb = $excel.ActiveWorkbook
$dsr = $wb.Worksheets.Item("DSR")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# DSR sheet: rewording of the review comment in D8
$dsr.Range("D8").Value = "How to get it?"

# Sheet1: new "Payment Mode" reference list in column J (values first, header last,
# to mirror the order the new strings were authored in)
$sheet1.Range("J13").Value = "Easy Diner"
$sheet1.Range("J14").Value = "Gpay"
$sheet1.Range("J9").Value = "dotpe googlepay"
$sheet1.Range("J6").Value = "Dineout"
$sheet1.Range("J22").Value = "UPI Payment"
$sheet1.Range("J10").Value = "dotpe paytm"
$sheet1.Range("J17").Value = "Paytm Deals"
$sheet1.Range("J23").Value = "Zomato"
$sheet1.Range("J3").Value = "BTC"
$sheet1.Range("J15").Value = "magicpin"
$sheet1.Range("J11").Value = "dotpe phonepe"
$sheet1.Range("J8").Value = "dotpe debitCard"
$sheet1.Range("J24").Value = "Zomato Pro"
$sheet1.Range("J2").Value = "BOOK MY SHOW"
$sheet1.Range("J20").Value = "SWIGGY Dineout"
$sheet1.Range("J4").Value = "Cash Sales"
$sheet1.Range("J12").Value = "dotpe simpl"
$sheet1.Range("J5").Value = "Credit Card Sales"
$sheet1.Range("J21").Value = "Thrive"
$sheet1.Range("J18").Value = "Sodexo"
$sheet1.Range("J16").Value = "PAYTM"
$sheet1.Range("J7").Value = "dotpe creditCard"
$sheet1.Range("J19").Value = "Swiggy"

# Header cell J1 (bold, matching B1/D1/F1/H1)
$sheet1.Range("J1").Value = "Payment Mode"
$sheet1.Range("J1").Font.Bold = $true
